$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking text (e.g. "1.020", "27.723.10") that must
# stay plain text, so force the whole data range to Text format before writing.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.689.49"
$ws.Range("E2").Value = "  +0.53%  "
$ws.Range("D3").Value = "1.855.96"
$ws.Range("E3").Value = "  +0.31%  "
$ws.Range("D4").Value = "1.018"
$ws.Range("E4").Value = "  -1.20%  "
$ws.Range("D5").Value = "320.02"
$ws.Range("E5").Value = "  -0.26%  "
$ws.Range("D6").Value = "1.017"
$ws.Range("E6").Value = "  -0.92%  "
$ws.Range("D7").Value = "0.4357"
$ws.Range("E7").Value = "  -0.62%  "
$ws.Range("D8").Value = "0.3776"
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("D9").Value = "0.07407"
$ws.Range("E9").Value = "  +0.15%  "
$ws.Range("D10").Value = "0.8811"
$ws.Range("E10").Value = "  +0.66%  "
$ws.Range("D11").Value = "21.59"
$ws.Range("E11").Value = "  +0.35%  "
$ws.Range("D12").Value = "1.875.77"
$ws.Range("E12").Value = "  +1.22%  "
$ws.Range("D13").Value = "5.482"
$ws.Range("E13").Value = "  -0.77%  "
$ws.Range("D14").Value = "6.729"
$ws.Range("E14").Value = "  +0.62%  "
$ws.Range("D15").Value = "0.07090"
$ws.Range("E15").Value = "  -1.63%  "
$ws.Range("D16").Value = "86.91"
$ws.Range("E16").Value = "  +4.87%  "
$ws.Range("E17").Value = "  -1.01%  "
$ws.Range("D18").Value = "0.000009047"
$ws.Range("E18").Value = "  +0.33%  "
$ws.Range("D19").Value = "1.017"
$ws.Range("E19").Value = "  -0.91%  "
$ws.Range("D20").Value = "15.38"
$ws.Range("E20").Value = "  -0.21%  "
$ws.Range("D21").Value = "27.704.52"
$ws.Range("E21").Value = "  +0.52%  "
$ws.Range("D22").Value = "5.284"
$ws.Range("E22").Value = "  +0.69%  "
$ws.Range("D23").Value = "11.12"
$ws.Range("E23").Value = "  -1.90%  "
$ws.Range("D24").Value = "2.090.51"
$ws.Range("E24").Value = "  +0.95%  "
$ws.Range("D25").Value = "2.032"
$ws.Range("E25").Value = "  +6.02%  "
$ws.Range("D26").Value = "156.84"
$ws.Range("E26").Value = "  -0.55%  "
$ws.Range("D27").Value = "18.67"
$ws.Range("E27").Value = "  -0.26%  "
$ws.Range("B28").Value = "InternetComputer(DFINITY)"
$ws.Range("C28").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D28").Value = "5.366"
$ws.Range("E28").Value = "  +1.84%  "
$ws.Range("B29").Value = "LidoDAOToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D29").Value = "1.993"
$ws.Range("E29").Value = "  +0.25%  "
$ws.Range("D30").Value = "120.79"
$ws.Range("E30").Value = "  +2.90%  "
$ws.Range("D31").Value = "0.09041"
$ws.Range("E31").Value = "  -0.11%  "
$ws.Range("D32").Value = "1.216"
$ws.Range("E32").Value = "  +1.66%  "
$ws.Range("D33").Value = "0.7690"
$ws.Range("E33").Value = "  +0.81%  "
$ws.Range("D34").Value = "3.027"
$ws.Range("E34").Value = "  +5.28%  "
$ws.Range("D35").Value = "4.548"
$ws.Range("E35").Value = "  +0.75%  "
$ws.Range("D36").Value = "1.018"
$ws.Range("E36").Value = "  -0.90%  "
$ws.Range("D37").Value = "1.139"
$ws.Range("E37").Value = "  -0.81%  "
$ws.Range("D38").Value = "0.01977"
$ws.Range("E38").Value = "  +0.05%  "
$ws.Range("D39").Value = "0.05295"
$ws.Range("E39").Value = "  +0.00%  "
$ws.Range("D40").Value = "2.870"
$ws.Range("E40").Value = "  +2.69%  "
$ws.Range("E41").Value = "  +0.78%  "
$ws.Range("D42").Value = "6.949"
$ws.Range("E42").Value = "  +3.03%  "
$ws.Range("D43").Value = "0.1677"
$ws.Range("E43").Value = "  +0.08%  "
$ws.Range("D44").Value = "8.686"
$ws.Range("E44").Value = "  +2.27%  "
$ws.Range("D45").Value = "109.98"
$ws.Range("E45").Value = "  +1.14%  "
$ws.Range("D46").Value = "10.74"
$ws.Range("E46").Value = "  +1.08%  "
$ws.Range("D47").Value = "1.709"
$ws.Range("E47").Value = "  +0.07%  "
$ws.Range("D48").Value = "0.4717"
$ws.Range("E48").Value = "  +1.47%  "
$ws.Range("D49").Value = "1.018"
$ws.Range("E49").Value = "  -1.08%  "
$ws.Range("D50").Value = "0.06486"
$ws.Range("E50").Value = "  +1.37%  "
$ws.Range("D51").Value = "1.852"
$ws.Range("E51").Value = "  +0.30%  "
